# Adapt column header formatting to respective input file names.
# - rename "<name>_old" columns to "<name>_FV2404"
# - rename "<name>_new" columns to "<name>_FV2410"
# - turn the data range into an Excel Table ("Table1")
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldSuffix = "_old"
$newSuffix = "_new"
$fv2404 = "_FV2404"
$fv2410 = "_FV2410"

# Header row is row 1; determine the number of used columns dynamically
$usedRange = $ws.UsedRange
$lastCol = $usedRange.Columns.Count

for ($c = 1; $c -le $lastCol; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if (($val -ne $null) -and ($val -is [string])) {
        if ($val.EndsWith($oldSuffix)) {
            $base = $val.Substring(0, $val.Length - $oldSuffix.Length)
            $cell.Value = $base + $fv2404
        } elseif ($val.EndsWith($newSuffix)) {
            $base = $val.Substring(0, $val.Length - $newSuffix.Length)
            $cell.Value = $base + $fv2410
        }
    }
}

# Determine the full used range for the table (header row + data rows)
$lastDataRow = $usedRange.Rows.Count
$tableRange = $ws.Range($ws.Cells.Item(1, 1), $ws.Cells.Item($lastDataRow, $lastCol))

# Create the Excel table (ListObject) over the data
$listObj = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$listObj.Name = "Table1"
$listObj.TableStyle = ""

# Freeze the header row (pane split below row 1)
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
